$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 113.69
$ws.Range("C7").Formula = "=0.06+118.62+0.9"
$ws.Range("C8").Value = 11.06
